# Fix: a detailed sub-group ("CX_SY_NO_HANDLER") was being skipped in the
# BusinessLogicException/Steplet pattern hierarchy, which caused the values
# of the skipped group to show up as the values of the next group instead.
# The fix inserts the missing row for CX_SY_NO_HANDLER (as a child/parent of
# BusinessLogicException, just like the existing ThreadInUse/Steplet rows)
# into the TopPatterns table, pushing the rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing "Steplet" row (row 4), shifting the
# Steplet row and everything below it down by one row.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the CX_SY_NO_HANDLER sub-pattern.
# Type/Group/ID (A4:C4) and Serverity (G4) stay blank, matching the other
# "detail" rows (e.g. ThreadInUse) that hang off the BusinessLogicException
# parent.
$ws.Range("D4").Value = "BusinessLogicException"
$ws.Range("E4").Value = "CX_SY_NO_HANDLER"
$ws.Range("F4").Value = "com.syclo.agentry.BusinessLogicException: (?P<steplet>.*)Steplet - (|(?P<_time>.*) - )CX_SY_NO_HANDLER: An exception of type '(?P<exceptionType>.*)' occured, that was not caught anywhere in the call hierarchy. It was not handled locally or declared using a raising cl"

# Grow the Table1_3 table so it includes the newly-inserted row.
$lo = $ws.ListObjects("Table1_3")
$lo.Resize($ws.Range("A1:G12"))
